$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 44165
$ws.Range("J2").Value = 720
$ws.Range("K2").Value = 1200
$ws.Range("L2").Value = 1200
$ws.Range("M2").Value = 1200
$ws.Range("P2").Value = 1200

# Row 3
$ws.Range("D3").Value = 44165
$ws.Range("J3").Value = 750
$ws.Range("K3").Value = 1000
$ws.Range("L3").Value = 1000
$ws.Range("M3").Value = 1000
$ws.Range("P3").Value = 1000

# Row 4
$ws.Range("D4").Value = 44168
$ws.Range("J4").Value = 1200

# Row 5
$ws.Range("D5").Value = 44168
$ws.Range("J5").Value = 850

# Row 6
$ws.Range("D6").Value = 44169
$ws.Range("J6").Value = 950

# Row 7
$ws.Range("D7").Value = 44169
$ws.Range("J7").Value = 800

# Row 8
$ws.Range("D8").Value = 44167
$ws.Range("J8").Value = 1430
$ws.Range("L8").Value = 1300
$ws.Range("M8").Value = 1248
$ws.Range("P8").Value = 1248

# Row 9
$ws.Range("D9").Value = 44167
$ws.Range("J9").Value = 350

# Row 10
$ws.Range("D10").Value = 44179
$ws.Range("J10").Value = 980
$ws.Range("K10").Value = 1200
$ws.Range("L10").Value = 1200
$ws.Range("M10").Value = 1200
$ws.Range("O10").Value = "Región Metropolitana"
$ws.Range("P10").Value = 1200

# Row 11
$ws.Range("D11").Value = 44175
$ws.Range("I11").Value = "Primera"
$ws.Range("J11").Value = 1500
$ws.Range("K11").Value = 1300
$ws.Range("L11").Value = 1300
$ws.Range("M11").Value = 1300
$ws.Range("P11").Value = 1300

# Row 12
$ws.Range("D12").Value = 44175
$ws.Range("I12").Value = "Segunda"
$ws.Range("J12").Value = 1450
$ws.Range("K12").Value = 1000
$ws.Range("L12").Value = 1000
$ws.Range("M12").Value = 1000
$ws.Range("P12").Value = 1000

# Row 13
$ws.Range("D13").Value = 44159
$ws.Range("I13").Value = "Primera"
$ws.Range("J13").Value = 1100
$ws.Range("K13").Value = 1300
$ws.Range("L13").Value = 1300
$ws.Range("M13").Value = 1300
$ws.Range("P13").Value = 1300

# Row 14
$ws.Range("D14").Value = 44159
$ws.Range("I14").Value = "Segunda"
$ws.Range("J14").Value = 800
$ws.Range("K14").Value = 1000
$ws.Range("L14").Value = 1000
$ws.Range("M14").Value = 1000
$ws.Range("O14").Value = "Provincia de Quillota"
$ws.Range("P14").Value = 1000

# Row 15
$ws.Range("D15").Value = 44176
$ws.Range("J15").Value = 2500
$ws.Range("L15").Value = 1300
$ws.Range("M15").Value = 1256
$ws.Range("P15").Value = 1256

# Row 16
$ws.Range("D16").Value = 44176
$ws.Range("J16").Value = 1500

# Row 17
$ws.Range("D17").Value = 44172
$ws.Range("J17").Value = 600

# Row 18
$ws.Range("D18").Value = 44172
$ws.Range("J18").Value = 550

# Row 21
$ws.Range("D21").Value = 44181
$ws.Range("J21").Value = 1000
$ws.Range("K21").Value = 1300
$ws.Range("M21").Value = 1300
$ws.Range("P21").Value = 1300

# Row 22
$ws.Range("D22").Value = 44181
$ws.Range("J22").Value = 900
$ws.Range("K22").Value = 900
$ws.Range("L22").Value = 900
$ws.Range("M22").Value = 900
$ws.Range("P22").Value = 900

# Row 23
$ws.Range("D23").Value = 44160
$ws.Range("J23").Value = 750

# Row 24
$ws.Range("D24").Value = 44160
$ws.Range("J24").Value = 850

# Row 25
$ws.Range("D25").Value = 44162
$ws.Range("J25").Value = 1200
$ws.Range("K25").Value = 1300
$ws.Range("M25").Value = 1300
$ws.Range("P25").Value = 1300

# Row 26
$ws.Range("D26").Value = 44162
$ws.Range("J26").Value = 800

# Row 27
$ws.Range("D27").Value = 44174
$ws.Range("J27").Value = 2800
$ws.Range("K27").Value = 1200
$ws.Range("L27").Value = 1250
$ws.Range("M27").Value = 1221
$ws.Range("P27").Value = 1221

# Row 28
$ws.Range("D28").Value = 44174
$ws.Range("J28").Value = 1300

